$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column H (8th column), shifting PerturbationDE and
# everything after it one column to the right.
$ws.Columns.Item(8).Insert()

# New header / value for the inserted "ParetoSize" column.
$ws.Cells.Item(1, 8).Value = "ParetoSize"
$ws.Cells.Item(2, 8).Value = 20
$ws.Columns.Item(8).ColumnWidth = 9.43

# NumGenGA value changed from 50 to 75.
$ws.Cells.Item(2, 6).Value = 75

# NumProcessors (column M after the insert) gets a value of 50.
$ws.Cells.Item(2, 13).Value = 50

$ws.Range("M17").Select()
